$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: confusion_matrix (text values containing embedded newlines)
$ws.Range("B2").Value = "[[30288  2619]`n [  909   397]]"
$ws.Range("C2").Value = "[[30328  2579]`n [  878   428]]"
$ws.Range("D2").Value = "[[30548  2359]`n [  958   348]]"
$ws.Range("E2").Value = "[[30334  2573]`n [  720   586]]"

# Row 3: accuracy_score
$ws.Range("B3").Value = 0.8968813024289013
$ws.Range("C3").Value = 0.8989565369888639
$ws.Range("D3").Value = 0.9030485487972408
$ws.Range("E3").Value = 0.9037500365358198

# Row 4: f1_score
$ws.Range("B4").Value = 0.1837112447940768
$ws.Range("C4").Value = 0.1984697426385347
$ws.Range("D4").Value = 0.1734363319212559
$ws.Range("E4").Value = 0.2624860022396416
